$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from H1 into I1:J1 before setting values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels
$ws.Range("I1").Value = 'DS_ESTADO_CIVIL'
$ws.Range("J1").Value = 'DS_GRAU_INSTRUCAO'

# Data values for DS_ESTADO_CIVIL (I) and DS_GRAU_INSTRUCAO (J), rows 2-84
$estadoCivil = @(
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'SOLTEIRO(A)',
    'DIVORCIADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'DIVORCIADO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'VIÚVO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'DIVORCIADO(A)',
    'DIVORCIADO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'DIVORCIADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'DIVORCIADO(A)',
    'CASADO(A)',
    'DIVORCIADO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'CASADO(A)',
    'SOLTEIRO(A)'
)

$grauInstrucao = @(
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO FUNDAMENTAL INCOMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO FUNDAMENTAL INCOMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR INCOMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR INCOMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'ENSINO FUNDAMENTAL INCOMPLETO',
    'ENSINO MÉDIO INCOMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO FUNDAMENTAL INCOMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR INCOMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO FUNDAMENTAL INCOMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR INCOMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR INCOMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR INCOMPLETO',
    'SUPERIOR INCOMPLETO',
    'ENSINO MÉDIO INCOMPLETO',
    'ENSINO FUNDAMENTAL INCOMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR COMPLETO',
    'SUPERIOR INCOMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO FUNDAMENTAL COMPLETO',
    'SUPERIOR COMPLETO',
    'ENSINO MÉDIO COMPLETO',
    'SUPERIOR COMPLETO'
)

for ($i = 0; $i -lt $estadoCivil.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $estadoCivil[$i]
    $ws.Cells.Item($row, 10).Value = $grauInstrucao[$i]
}